$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.100659
$ws.Range("H2").Value = 30.301977
$ws.Range("I2").Value = 0.3328245842863797
$ws.Range("J2").Value = 0.3328245842863797
$ws.Range("M2").Value = 4.372354666666666
$ws.Range("N2").Value = 13.117064
$ws.Range("O2").Value = 0.04826584573009855
$ws.Range("P2").Value = 0.04826584573009855
$ws.Range("Q2").Value = 44.16366351505867
$ws.Range("R2").Value = 397.472971635528
$ws.Range("S2").Value = 0.01606406004035058
$ws.Range("T2").Value = 0.01606406004035058
$ws.Range("G3").Value = 10.100659
$ws.Range("H3").Value = 30.301977
$ws.Range("I3").Value = 0.3328245842863797
$ws.Range("J3").Value = 0.3328245842863797
$ws.Range("O3").Value = 0.0323794466022692
$ws.Range("P3").Value = 0.0323794466022692
$ws.Range("Q3").Value = 29.62747182641167
$ws.Range("R3").Value = 266.6472464377051
$ws.Range("S3").Value = 0.01077667585482328
$ws.Range("T3").Value = 0.01077667585482328
$ws.Range("G4").Value = 10.100659
$ws.Range("H4").Value = 30.301977
$ws.Range("I4").Value = 0.3328245842863797
$ws.Range("J4").Value = 0.3328245842863797
$ws.Range("M4").Value = 71.62252933333333
$ws.Range("N4").Value = 214.867588
$ws.Range("O4").Value = 0.7906316424777964
$ws.Range("P4").Value = 0.7906316424777964
$ws.Range("Q4").Value = 723.4347455134973
$ws.Range("R4").Value = 6510.912709621476
$ws.Range("S4").Value = 0.2631416477313301
$ws.Range("T4").Value = 0.2631416477313301
$ws.Range("G5").Value = 10.100659
$ws.Range("H5").Value = 30.301977
$ws.Range("I5").Value = 0.3328245842863797
$ws.Range("J5").Value = 0.3328245842863797
$ws.Range("M5").Value = 11.66089366666667
$ws.Range("N5").Value = 34.982681
$ws.Range("O5").Value = 0.1287230651898359
$ws.Range("P5").Value = 0.1287230651898359
$ws.Range("Q5").Value = 117.7827105622597
$ws.Range("R5").Value = 1060.044395060337
$ws.Range("S5").Value = 0.04284220065987569
$ws.Range("T5").Value = 0.04284220065987569
$ws.Range("I6").Value = 0.4180918757349671
$ws.Range("J6").Value = 0.4180918757349671
$ws.Range("M6").Value = 4.372354666666666
$ws.Range("N6").Value = 13.117064
$ws.Range("O6").Value = 0.04826584573009855
$ws.Range("P6").Value = 0.04826584573009855
$ws.Range("Q6").Value = 55.47808001602133
$ws.Range("R6").Value = 499.302720144192
$ws.Range("S6").Value = 0.02017955797523145
$ws.Range("T6").Value = 0.02017955797523145
$ws.Range("I7").Value = 0.4180918757349671
$ws.Range("J7").Value = 0.4180918757349671
$ws.Range("O7").Value = 0.0323794466022692
$ws.Range("P7").Value = 0.0323794466022692
$ws.Range("S7").Value = 0.01353758356520294
$ws.Range("T7").Value = 0.01353758356520294
$ws.Range("I8").Value = 0.4180918757349671
$ws.Range("J8").Value = 0.4180918757349671
$ws.Range("M8").Value = 71.62252933333333
$ws.Range("N8").Value = 214.867588
$ws.Range("O8").Value = 0.7906316424777964
$ws.Range("P8").Value = 0.7906316424777964
$ws.Range("Q8").Value = 908.7735822523626
$ws.Range("R8").Value = 8178.962240271265
$ws.Range("S8").Value = 0.3305566664189598
$ws.Range("T8").Value = 0.3305566664189598
$ws.Range("I9").Value = 0.4180918757349671
$ws.Range("J9").Value = 0.4180918757349671
$ws.Range("M9").Value = 11.66089366666667
$ws.Range("N9").Value = 34.982681
$ws.Range("O9").Value = 0.1287230651898359
$ws.Range("P9").Value = 0.1287230651898359
$ws.Range("Q9").Value = 147.9578033386853
$ws.Range("R9").Value = 1331.620230048168
$ws.Range("S9").Value = 0.05381806777557294
$ws.Range("T9").Value = 0.05381806777557294
$ws.Range("G10").Value = 4.721016333333334
$ws.Range("H10").Value = 14.163049
$ws.Range("I10").Value = 0.1555611667071302
$ws.Range("J10").Value = 0.1555611667071302
$ws.Range("M10").Value = 4.372354666666666
$ws.Range("N10").Value = 13.117064
$ws.Range("O10").Value = 0.04826584573009855
$ws.Range("P10").Value = 0.04826584573009855
$ws.Range("Q10").Value = 20.64195779645956
$ws.Range("R10").Value = 185.777620168136
$ws.Range("S10").Value = 0.007508291273880491
$ws.Range("T10").Value = 0.007508291273880489
$ws.Range("G11").Value = 4.721016333333334
$ws.Range("H11").Value = 14.163049
$ws.Range("I11").Value = 0.1555611667071302
$ws.Range("J11").Value = 0.1555611667071302
$ws.Range("O11").Value = 0.0323794466022692
$ws.Range("P11").Value = 0.0323794466022692
$ws.Range("Q11").Value = 13.84778739762056
$ws.Range("R11").Value = 124.630086578585
$ws.Range("S11").Value = 0.005036984490780221
$ws.Range("T11").Value = 0.00503698449078022
$ws.Range("G12").Value = 4.721016333333334
$ws.Range("H12").Value = 14.163049
$ws.Range("I12").Value = 0.1555611667071302
$ws.Range("J12").Value = 0.1555611667071302
$ws.Range("M12").Value = 71.62252933333333
$ws.Range("N12").Value = 214.867588
$ws.Range("O12").Value = 0.7906316424777964
$ws.Range("P12").Value = 0.7906316424777964
$ws.Range("Q12").Value = 338.1311308173125
$ws.Range("R12").Value = 3043.180177355812
$ws.Range("S12").Value = 0.1229915807394207
$ws.Range("T12").Value = 0.1229915807394207
$ws.Range("G13").Value = 4.721016333333334
$ws.Range("H13").Value = 14.163049
$ws.Range("I13").Value = 0.1555611667071302
$ws.Range("J13").Value = 0.1555611667071302
$ws.Range("M13").Value = 11.66089366666667
$ws.Range("N13").Value = 34.982681
$ws.Range("O13").Value = 0.1287230651898359
$ws.Range("P13").Value = 0.1287230651898359
$ws.Range("Q13").Value = 55.05126946159656
$ws.Range("R13").Value = 495.461425154369
$ws.Range("S13").Value = 0.02002431020304886
$ws.Range("T13").Value = 0.02002431020304885
$ws.Range("G14").Value = 2.838244666666667
$ws.Range("H14").Value = 8.514734000000001
$ws.Range("I14").Value = 0.09352237327152295
$ws.Range("J14").Value = 0.09352237327152294
$ws.Range("M14").Value = 4.372354666666666
$ws.Range("N14").Value = 13.117064
$ws.Range("O14").Value = 0.04826584573009855
$ws.Range("P14").Value = 0.04826584573009855
$ws.Range("Q14").Value = 12.40981231344178
$ws.Range("R14").Value = 111.688310820976
$ws.Range("S14").Value = 0.004513936440636019
$ws.Range("T14").Value = 0.004513936440636018
$ws.Range("G15").Value = 2.838244666666667
$ws.Range("H15").Value = 8.514734000000001
$ws.Range("I15").Value = 0.09352237327152295
$ws.Range("J15").Value = 0.09352237327152294
$ws.Range("O15").Value = 0.0323794466022692
$ws.Range("P15").Value = 0.0323794466022692
$ws.Range("Q15").Value = 8.325200751567779
$ws.Range("R15").Value = 74.92680676411001
$ws.Range("S15").Value = 0.003028202691462765
$ws.Range("T15").Value = 0.003028202691462765
$ws.Range("G16").Value = 2.838244666666667
$ws.Range("H16").Value = 8.514734000000001
$ws.Range("I16").Value = 0.09352237327152295
$ws.Range("J16").Value = 0.09352237327152294
$ws.Range("M16").Value = 71.62252933333333
$ws.Range("N16").Value = 214.867588
$ws.Range("O16").Value = 0.7906316424777964
$ws.Range("P16").Value = 0.7906316424777964
$ws.Range("Q16").Value = 203.2822618935102
$ws.Range("R16").Value = 1829.540357041592
$ws.Range("S16").Value = 0.07394174758808575
$ws.Range("T16").Value = 0.07394174758808573
$ws.Range("G17").Value = 2.838244666666667
$ws.Range("H17").Value = 8.514734000000001
$ws.Range("I17").Value = 0.09352237327152295
$ws.Range("J17").Value = 0.09352237327152294
$ws.Range("M17").Value = 11.66089366666667
$ws.Range("N17").Value = 34.982681
$ws.Range("O17").Value = 0.1287230651898359
$ws.Range("P17").Value = 0.1287230651898359
$ws.Range("Q17").Value = 33.09646925798378
$ws.Range("R17").Value = 297.868223321854
$ws.Range("S17").Value = 0.01203848655133842
$ws.Range("T17").Value = 0.01203848655133841

Write-Output "Updated 174 cells"